# Generate Report for Handoff
# Rewrites the localization-status workbook so the tracked file is now
# ad230e1d-87a4-4a71-b460-5e004e553def.md (was 0f51c9b3-f0a2-4be6-b663-558a8ba7405b.md)
# and refreshes the handoff timestamps / handback state for a fresh
# "Ready for handoff" run (target files / handback files / handback dates
# are cleared back out since nothing has come back yet).

$wb = $excel.ActiveWorkbook

$oldUuid = "0f51c9b3-f0a2-4be6-b663-558a8ba7405b"
$newUuid = "ad230e1d-87a4-4a71-b460-5e004e553def"
$oldHash = "b53ab6119df9594e81535c4c4b1de34c79abd09a"
$newHash = "430d3ba8d07312508c71c3cdcea38a19f6bc80e8"

$overviewHyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1811585a046efea986da45231734b31b65e72af5/e2e/$oldUuid.md"
$zhcnHyperlinkAddr      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1811585a046efea986da45231734b31b65e72af5/e2e/$oldUuid.md"
$dedeHyperlinkAddr      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1811585a046efea986da45231734b31b65e72af5/e2e/$oldUuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid.md"
$wsOverview.Range("G2").Value = "2016-08-30 19:14:40"

# Refresh the B2 hyperlink's display text (target/address is unchanged).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewHyperlinkAddr, "", "", "e2e\$newUuid.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid.md"
$wsZhCn.Range("G2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 19:14:36"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Only the A2 hyperlink survives (the I2 "Latest Target File" hyperlink is
# gone now that the column is blank again); refresh A2's display text.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnHyperlinkAddr, "", "", "$newUuid.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid.md"
$wsDeDe.Range("G2").Value = "$newUuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-30 19:14:40"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeHyperlinkAddr, "", "", "$newUuid.md")

# ---------------------------------------------------------------------
# Column widths on zh-cn / de-de: "Latest Target File" (I) and
# "Latest Handback File" (J) shrink back to their auto-fit width now that
# their sample values are gone.
# ---------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(9).ColumnWidth = 18.6506053379604
    $ws.Columns.Item(10).ColumnWidth = 21.7054770333426
}
